$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$changes = @(
    @("D2", "254.16"),
    @("E2", "3.14%"),
    @("G2", "6"),
    @("D3", "27.91"),
    @("E3", "-5.14%"),
    @("G3", "6"),
    @("D4", "5.332"),
    @("E4", "3.57%"),
    @("G4", "6"),
    @("D5", "0.05842"),
    @("E5", "0.86%"),
    @("G5", "6"),
    @("D6", "6.710"),
    @("E6", "1.04%"),
    @("G6", "6"),
    @("D7", "0.8687"),
    @("E7", "1.73%"),
    @("G7", "6"),
    @("D8", "0.9161"),
    @("E8", "5.73%"),
    @("G8", "6"),
    @("D9", "0.1429"),
    @("E9", "3.91%"),
    @("G9", "6"),
    @("D10", "0.07172"),
    @("E10", "1.53%"),
    @("G10", "6"),
    @("D11", "0.03238"),
    @("E11", "-0.87%"),
    @("G11", "6"),
    @("D12", "0.09240"),
    @("E12", "-1.40%"),
    @("G12", "6"),
    @("E13", "1.14%"),
    @("G13", "6"),
    @("D14", "0.0006068"),
    @("E14", "1.60%"),
    @("G14", "6"),
    @("D15", "0.006045"),
    @("E15", "-1.02%"),
    @("G15", "6"),
    @("E16", "0.33%"),
    @("G16", "6"),
    @("E17", "1.44%"),
    @("G17", "6"),
    @("D18", "2.251"),
    @("E18", "1.86%"),
    @("G18", "6"),
    @("E19", "-1.01%"),
    @("G19", "6"),
    @("D20", "0.03447"),
    @("E20", "3.15%"),
    @("G20", "6"),
    @("D21", "0.1302"),
    @("E21", "1.62%"),
    @("G21", "6"),
    @("D22", "3.527"),
    @("E22", "11.45%"),
    @("G22", "6"),
    @("D23", "0.04143"),
    @("E23", "0.11%"),
    @("G23", "6"),
    @("D24", "0.1345"),
    @("E24", "-3.84%"),
    @("G24", "6"),
    @("D25", "0.005107"),
    @("E25", "23.30%"),
    @("G25", "6"),
    @("D26", "0.001225"),
    @("E26", "-0.12%"),
    @("G26", "6"),
    @("E27", "-0.71%"),
    @("G27", "6"),
    @("D28", "0.0001937"),
    @("E28", "34.22%"),
    @("G28", "6"),
    @("G29", "6"),
    @("G30", "6"),
    @("G31", "6"),
    @("G32", "6"),
    @("G33", "6"),
    @("G34", "6"),
    @("G35", "6"),
    @("G36", "6"),
    @("G37", "6"),
    @("G38", "6"),
    @("G39", "6"),
    @("D40", "0.03848"),
    @("E40", "2.90%"),
    @("G40", "6"),
    @("D41", "0.1098"),
    @("E41", "2.53%"),
    @("G41", "6"),
    @("D42", "0.002390"),
    @("E42", "8.76%"),
    @("G42", "6"),
    @("D43", "0.002950"),
    @("E43", "-48.82%"),
    @("G43", "6"),
    @("D44", "0.009964"),
    @("E44", "8.61%"),
    @("G44", "6"),
    @("D45", "0.00005287"),
    @("E45", "0.25%"),
    @("G45", "6"),
    @("G46", "6"),
    @("D47", "0.09997"),
    @("E47", "72.57%"),
    @("G47", "6"),
    @("D48", "0.002185"),
    @("E48", "0.55%"),
    @("G48", "6"),
    @("G49", "6"),
    @("D50", "0.0001999"),
    @("G50", "6"),
    @("G51", "6"),
)

foreach ($item in $changes) {
    $cellRef = $item[0]
    $val = $item[1]
    $c = $ws.Range($cellRef)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.ClearFormats()
}